$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These "Price" cells receive plain decimal numbers as their new value
# (e.g. "18.33"). Left alone, Excel auto-parses a bare decimal into a
# numeric literal (losing trailing zeros / exact text, e.g. "0.0790" ->
# 0.079). The source feed stores Price as literal text, so force a Text
# number format on just those cells before assigning, preserving the
# exact string content.
$textForcedCells = @(
    "D5",
    "D10",
    "D11",
    "D15",
    "D18",
    "D20",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D28",
    "D29",
    "D31",
    "D33",
    "D35",
    "D37",
    "D39",
    "D41",
    "D44",
    "D46",
    "D47",
    "D48",
    "D51"
)
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.028.47'
$ws.Range('E2').Value = '  -0.05%  '
$ws.Range('D3').Value = '1.632.78'
$ws.Range('E3').Value = '  -0.70%  '
$ws.Range('E4').Value = '  +0.37%  '
$ws.Range('D5').Value = '214.59'
$ws.Range('E5').Value = '  -0.87%  '
$ws.Range('E6').Value = '  -1.05%  '
$ws.Range('E7').Value = '  +0.36%  '
$ws.Range('E8').Value = '  -2.54%  '
$ws.Range('E9').Value = '  -3.21%  '
$ws.Range('D10').Value = '18.33'
$ws.Range('E10').Value = '  -6.78%  '
$ws.Range('D11').Value = '0.0790'
$ws.Range('E11').Value = '  -0.69%  '
$ws.Range('D12').Value = '1.859.85'
$ws.Range('E12').Value = '  -0.67%  '
$ws.Range('D13').Value = '1.629.62'
$ws.Range('E13').Value = '  -2.04%  '
$ws.Range('E14').Value = '  -2.74%  '
$ws.Range('D15').Value = '0.524'
$ws.Range('E15').Value = '  -3.81%  '
$ws.Range('D16').Value = '25.991.84'
$ws.Range('E16').Value = '  -0.23%  '
$ws.Range('D17').Value = '0.0₃0741'
$ws.Range('E17').Value = '  -3.07%  '
$ws.Range('D18').Value = '61.38'
$ws.Range('E19').Value = '  +0.34%  '
$ws.Range('D20').Value = '190.44'
$ws.Range('E20').Value = '  -2.61%  '
$ws.Range('E21').Value = '  -2.27%  '
$ws.Range('D22').Value = '9.62'
$ws.Range('E22').Value = '  -3.28%  '
$ws.Range('D23').Value = '6.09'
$ws.Range('E23').Value = '  -2.52%  '
$ws.Range('D24').Value = '0.133'
$ws.Range('E24').Value = '  +0.51%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '143.75'
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = '1.79'
$ws.Range('E26').Value = '  -0.98%  '
$ws.Range('E27').Value = '  +0.20%  '
$ws.Range('D28').Value = '6.74'
$ws.Range('E28').Value = '  -2.34%  '
$ws.Range('D29').Value = '15.14'
$ws.Range('E29').Value = '  -2.70%  '
$ws.Range('E30').Value = '  -1.38%  '
$ws.Range('D31').Value = '0.0482'
$ws.Range('E31').Value = '  -3.30%  '
$ws.Range('E32').Value = '  -4.05%  '
$ws.Range('D33').Value = '3.11'
$ws.Range('E33').Value = '  -5.45%  '
$ws.Range('E34').Value = '  -2.30%  '
$ws.Range('D35').Value = '1.48'
$ws.Range('E35').Value = '  -3.79%  '
$ws.Range('D36').Value = '1.132.99'
$ws.Range('E36').Value = '  -0.10%  '
$ws.Range('D37').Value = '0.856'
$ws.Range('E37').Value = '  -5.60%  '
$ws.Range('E38').Value = '  -0.90%  '
$ws.Range('D39').Value = '0.517'
$ws.Range('E39').Value = '  -4.57%  '
$ws.Range('E40').Value = '  -1.65%  '
$ws.Range('D41').Value = '98.12'
$ws.Range('E41').Value = '  -1.32%  '
$ws.Range('D43').Value = '1.770.21'
$ws.Range('E43').Value = '  -0.63%  '
$ws.Range('D44').Value = '5.20'
$ws.Range('E44').Value = '  -5.43%  '
$ws.Range('E45').Value = '  -2.36%  '
$ws.Range('D46').Value = '54.72'
$ws.Range('E46').Value = '  -3.61%  '
$ws.Range('D47').Value = '0.0525'
$ws.Range('E47').Value = '  +0.09%  '
$ws.Range('D48').Value = '1.48'
$ws.Range('E48').Value = '  +0.54%  '
$ws.Range('E49').Value = '  +0.12%  '
$ws.Range('E50').Value = '  +0.41%  '
$ws.Range('D51').Value = '7.46'
$ws.Range('E51').Value = '  -3.70%  '
